# semana 32 de 2025
# Adds a new weekly column "AI" (week 32) to the IRA hospitalario weekly
# tracking sheet, mirroring the header style used by the other week
# columns, and fills in the reported counts for the UPGDs that reported
# during week 32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell AI1: new week number "32" ------------------------------
# The other week-number headers (D1:AH1) are stored as bold, centered text.
# Use an apostrophe-prefixed value so the engine stores "32" as text
# (matching the existing header cells) instead of a number, then re-apply
# the bold + centered formatting used throughout row 1.
$ws.Range("AI1").Value = "'32"
$ws.Range("AI1").Font.Bold = $true
$ws.Range("AI1").HorizontalAlignment = -4108

# --- Data cells: reported counts for week 32 -----------------------------
$ws.Range("AI2").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AI6").Value = 30
$ws.Range("AI7").Value = 1
$ws.Range("AI8").Value = 25
$ws.Range("AI9").Value = 0
$ws.Range("AI10").Value = 0
$ws.Range("AI13").Value = 0
$ws.Range("AI15").Value = 0
$ws.Range("AI16").Value = 0
$ws.Range("AI17").Value = 0
$ws.Range("AI23").Value = 0
$ws.Range("AI25").Value = 1
$ws.Range("AI27").Value = 0
$ws.Range("AI28").Value = 36
$ws.Range("AI29").Value = 3
$ws.Range("AI30").Value = 19
$ws.Range("AI31").Value = 0
$ws.Range("AI32").Value = 0
$ws.Range("AI34").Value = 3
$ws.Range("AI35").Value = 18
$ws.Range("AI36").Value = 0
$ws.Range("AI37").Value = 0
$ws.Range("AI38").Value = 0
$ws.Range("AI40").Value = 0
$ws.Range("AI41").Value = 0
$ws.Range("AI42").Value = 0
$ws.Range("AI43").Value = 0
$ws.Range("AI45").Value = 0
$ws.Range("AI46").Value = 0
$ws.Range("AI47").Value = 0
$ws.Range("AI48").Value = 0
$ws.Range("AI49").Value = 0
$ws.Range("AI50").Value = 0
$ws.Range("AI51").Value = 0
$ws.Range("AI53").Value = 0
$ws.Range("AI54").Value = 0
$ws.Range("AI55").Value = 0
$ws.Range("AI56").Value = 0
$ws.Range("AI57").Value = 0
$ws.Range("AI58").Value = 0
